$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "28.019.30"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -1.99%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.831.47"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -1.01%  "

$ws.Range("E4").Value = "  -0.13%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "324.28"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -2.74%  "

$ws.Range("E6").Value = "  -0.06%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4662"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.03%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3866"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -1.42%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.07883"
$c.Style = "Normal"

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.9602"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -2.52%  "

$ws.Range("E11").Value = "  -1.61%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "1.786.66"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -10.78%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "5.671"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -3.11%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "6.907"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -1.71%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.06836"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +0.16%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "87.37"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -0.36%  "

$ws.Range("E17").Value = "  -0.14%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.000009929"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -1.38%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "16.58"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -2.66%  "

$ws.Range("E20").Value = "  +0.02%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "28.011.78"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -2.10%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.318"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -1.49%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "10.96"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -2.33%  "

$ws.Range("E24").Value = "  -1.96%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.071.32"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -7.32%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "154.01"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +0.25%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "19.09"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -1.40%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "5.724"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -6.30%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.961"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -2.64%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "117.53"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +0.07%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.09268"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -1.71%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.9349"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -4.16%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "5.282"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -1.59%  "

$ws.Range("E34").Value = "  -2.03%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "3.296"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -5.89%  "

$ws.Range("E36").Value = "  -3.98%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.02141"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -2.49%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "1.142"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -1.91%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "7.784"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +2.47%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.5579"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -2.19%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "9.855"
$c.Style = "Normal"

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.1757"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -1.89%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "11.57"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -2.48%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.5257"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -2.31%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.07010"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -2.10%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "2.127"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -11.06%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "1.824"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -4.26%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "113.22"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -0.02%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.101"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -12.12%  "

$ws.Range("E50").Value = "  -0.04%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "2.323"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.64%  "
